$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")
Write-Host ("A1 value: " + $ws.Range("B1").Value)
